# Atualização automática de preços de eletricidade
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Date (A2) - serial date value
$ws.Range("A2").Value = 45938

# Hourly prices B2:Z2
$ws.Range("B2").Value = 111.56
$ws.Range("C2").Value = 104.25
$ws.Range("D2").Value = 101.15
$ws.Range("E2").Value = 99.75
$ws.Range("F2").Value = 100
$ws.Range("G2").Value = 97.56
$ws.Range("H2").Value = 105.24
$ws.Range("I2").Value = 129.15
$ws.Range("J2").Value = 132.25
$ws.Range("K2").Value = 109.85
$ws.Range("L2").Value = 98.03
$ws.Range("M2").Value = 79.73999999999999
$ws.Range("N2").Value = 76.19
$ws.Range("O2").Value = 68.2
$ws.Range("P2").Value = 63.23
$ws.Range("Q2").Value = 61.95
$ws.Range("R2").Value = 55.55
$ws.Range("S2").Value = 77.09999999999999
$ws.Range("T2").Value = 91.03
$ws.Range("U2").Value = 141.79
$ws.Range("V2").Value = 139.05
$ws.Range("W2").Value = 113.27
$ws.Range("X2").Value = 106.37
$ws.Range("Y2").Value = 96.87
$ws.Range("Z2").Value = 98.3

# Slot summary columns (AA-AG)
# AA2 (Slot_4h_max) unchanged: "20h-24h"
$ws.Range("AB2").Value = 113.89
# AC2 (Slot_2h_frist) unchanged: "20h-22h"
$ws.Range("AD2").Value = 126.16
$ws.Range("AE2").Value = "8h-10h"
$ws.Range("AF2").Value = 121.05
$ws.Range("AG2").Value = "5h-23h"
